# Apply "Added 12765 and 12785 surrogates" edit.
# The four existing surrogate rows (2-5) are cyclically rotated down by one
# row (row 2's data -> row 3, row 3 -> row 4, row 4 -> row 5, row 5 -> row 2),
# making room at the top of the block for the newly added surrogates while
# keeping the previously-present compounds (n-C12, iso-octane, 1,3,5-tmb,
# n-propylbenzene) in the sheet, reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 <- old Row 5 (n-propylbenzene)
$ws.Range("A2").Value = "n-propylbenzene"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("L2").Value = 5
$ws.Range("O2").Value = 1
$ws.Range("CP2").Value = 1

# Row 3 <- old Row 2 (n-C12)
$ws.Range("A3").Value = "n-C12"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 0
$ws.Range("CB3").Value = 0

# Row 4 <- old Row 3 (iso-octane)
$ws.Range("A4").Value = "iso-octane"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 1
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("CB4").Value = 1

# Row 5 <- old Row 4 (1,3,5-tmb)
$ws.Range("A5").Value = "1,3,5-tmb"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("L5").Value = 3
$ws.Range("N5").Value = 3
$ws.Range("O5").Value = 0
$ws.Range("CP5").Value = 0

# Leave the selection on the freshly edited row, matching the editor's
# final on-screen selection (entire row 2).
$ws.Range("A2:XFD2").Select()
